# Update the date line.
$d = $word.ActiveDocument
$d.Content.Find.Execute("2024-04-04 Thursday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2024-04-05 Friday", 2)

# Update the division-practice table cells in place (by row/column) so that
# duplicate old/new text values cannot cause a mismatched global replace.
$t = $d.Tables.Item(1)

$updates = @(
    @{Row=1;  Col=1; Text="96÷4=24, 0"},
    @{Row=1;  Col=2; Text="89÷9=9, 8"},
    @{Row=1;  Col=3; Text="17÷4=4, 1"},
    @{Row=1;  Col=4; Text="62÷6=10, 2"},
    @{Row=1;  Col=5; Text="68÷9=7, 5"},

    @{Row=5;  Col=1; Text="84÷7=12, 0"},
    @{Row=5;  Col=2; Text="42÷5=8, 2"},
    @{Row=5;  Col=3; Text="29÷4=7, 1"},
    @{Row=5;  Col=4; Text="18÷9=2, 0"},
    @{Row=5;  Col=5; Text="90÷8=11, 2"},

    @{Row=9;  Col=1; Text="21÷2=10, 1"},
    @{Row=9;  Col=2; Text="33÷3=11, 0"},
    @{Row=9;  Col=3; Text="56÷3=18, 2"},
    @{Row=9;  Col=4; Text="40÷9=4, 4"},
    @{Row=9;  Col=5; Text="13÷8=1, 5"},

    @{Row=13; Col=1; Text="88÷2=44, 0"},
    @{Row=13; Col=2; Text="75÷8=9, 3"},
    @{Row=13; Col=3; Text="87÷6=14, 3"},
    @{Row=13; Col=4; Text="66÷4=16, 2"},
    @{Row=13; Col=5; Text="42÷6=7, 0"},

    @{Row=17; Col=1; Text="64÷5=12, 4"},
    @{Row=17; Col=2; Text="53÷2=26, 1"},
    @{Row=17; Col=3; Text="25÷8=3, 1"},
    @{Row=17; Col=4; Text="45÷3=15, 0"},
    @{Row=17; Col=5; Text="74÷8=9, 2"}
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $cell.Range.Text = $u.Text
}
